$d = $word.ActiveDocument

$d.Content.Find.Execute("412÷8=51, 4", $true, $false, $false, $false, $false, $true, 1, $false, "638÷8=79, 6", 2)
$d.Content.Find.Execute("120÷6=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "369÷4=92, 1", 2)
$d.Content.Find.Execute("145÷2=72, 1", $true, $false, $false, $false, $false, $true, 1, $false, "747÷2=373, 1", 2)
$d.Content.Find.Execute("262÷2=131, 0", $true, $false, $false, $false, $false, $true, 1, $false, "582÷2=291, 0", 2)
$d.Content.Find.Execute("646÷4=161, 2", $true, $false, $false, $false, $false, $true, 1, $false, "561÷5=112, 1", 2)
$d.Content.Find.Execute("615÷4=153, 3", $true, $false, $false, $false, $false, $true, 1, $false, "505÷6=84, 1", 2)
$d.Content.Find.Execute("342÷4=85, 2", $true, $false, $false, $false, $false, $true, 1, $false, "544÷5=108, 4", 2)
$d.Content.Find.Execute("906÷6=151, 0", $true, $false, $false, $false, $false, $true, 1, $false, "331÷2=165, 1", 2)
$d.Content.Find.Execute("791÷8=98, 7", $true, $false, $false, $false, $false, $true, 1, $false, "793÷2=396, 1", 2)
$d.Content.Find.Execute("310÷7=44, 2", $true, $false, $false, $false, $false, $true, 1, $false, "604÷5=120, 4", 2)
$d.Content.Find.Execute("883÷6=147, 1", $true, $false, $false, $false, $false, $true, 1, $false, "623÷3=207, 2", 2)
$d.Content.Find.Execute("449÷8=56, 1", $true, $false, $false, $false, $false, $true, 1, $false, "604÷9=67, 1", 2)
$d.Content.Find.Execute("366÷5=73, 1", $true, $false, $false, $false, $false, $true, 1, $false, "849÷9=94, 3", 2)
$d.Content.Find.Execute("409÷9=45, 4", $true, $false, $false, $false, $false, $true, 1, $false, "631÷4=157, 3", 2)
$d.Content.Find.Execute("273÷8=34, 1", $true, $false, $false, $false, $false, $true, 1, $false, "984÷6=164, 0", 2)
$d.Content.Find.Execute("771÷2=385, 1", $true, $false, $false, $false, $false, $true, 1, $false, "425÷7=60, 5", 2)
$d.Content.Find.Execute("242÷8=30, 2", $true, $false, $false, $false, $false, $true, 1, $false, "543÷5=108, 3", 2)
$d.Content.Find.Execute("665÷6=110, 5", $true, $false, $false, $false, $false, $true, 1, $false, "385÷7=55, 0", 2)
$d.Content.Find.Execute("139÷9=15, 4", $true, $false, $false, $false, $false, $true, 1, $false, "421÷2=210, 1", 2)
$d.Content.Find.Execute("752÷8=94, 0", $true, $false, $false, $false, $false, $true, 1, $false, "754÷6=125, 4", 2)
$d.Content.Find.Execute("407÷8=50, 7", $true, $false, $false, $false, $false, $true, 1, $false, "265÷9=29, 4", 2)
$d.Content.Find.Execute("857÷4=214, 1", $true, $false, $false, $false, $false, $true, 1, $false, "509÷3=169, 2", 2)
$d.Content.Find.Execute("661÷6=110, 1", $true, $false, $false, $false, $false, $true, 1, $false, "486÷8=60, 6", 2)
$d.Content.Find.Execute("726÷2=363, 0", $true, $false, $false, $false, $false, $true, 1, $false, "301÷4=75, 1", 2)
$d.Content.Find.Execute("260÷8=32, 4", $true, $false, $false, $false, $false, $true, 1, $false, "910÷2=455, 0", 2)
